# Auto-generated edit script applying F-column ('想去人数') updates
# across all four worksheets, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3411
$ws.Range("F4").Value = 568
$ws.Range("F5").Value = 829
$ws.Range("F6").Value = 313
$ws.Range("F7").Value = 270
$ws.Range("F9").Value = 156
$ws.Range("F10").Value = 616
$ws.Range("F11").Value = 196
$ws.Range("F12").Value = 405
$ws.Range("F13").Value = 57
$ws.Range("F14").Value = 481
$ws.Range("F15").Value = 293
$ws.Range("F16").Value = 56
$ws.Range("F18").Value = 96
$ws.Range("F19").Value = 176

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 154
$ws.Range("F4").Value = 21
$ws.Range("F8").Value = 107
$ws.Range("F10").Value = 175
$ws.Range("F11").Value = 6
$ws.Range("F15").Value = 26
$ws.Range("F18").Value = 41

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6200
$ws.Range("F5").Value = 1769
$ws.Range("F6").Value = 108

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6200
$ws.Range("F5").Value = 1769
$ws.Range("F6").Value = 3411
$ws.Range("F7").Value = 154
$ws.Range("F8").Value = 108
$ws.Range("F10").Value = 568
$ws.Range("F11").Value = 829
$ws.Range("F12").Value = 313
$ws.Range("F13").Value = 270
$ws.Range("F16").Value = 21
$ws.Range("F18").Value = 156
$ws.Range("F21").Value = 617
$ws.Range("F22").Value = 107
$ws.Range("F23").Value = 196
$ws.Range("F25").Value = 405
$ws.Range("F26").Value = 175
$ws.Range("F27").Value = 57
$ws.Range("F28").Value = 481
$ws.Range("F29").Value = 6
$ws.Range("F30").Value = 293
$ws.Range("F31").Value = 56
$ws.Range("F35").Value = 96
$ws.Range("F37").Value = 26
$ws.Range("F40").Value = 41
$ws.Range("F41").Value = 176
